# Generate Report for handoff
#
# This script reproduces the localization-status report update:
#  - the previously-failed handoff file (c6618b96-...md) is renamed to a new
#    guid (f6f5a2a4-...md) and its status flips from "Handoff transform failed"
#    to "Ready for handoff"
#  - a brand-new file (ffff856d4578-...md) shows up, also "Ready for handoff"
#  - the .localization-config bookkeeping row moves down one row to make room
#  - the per-language detail sheets (zh-cn / de-de) gain "Latest Handoff File"
#    / "Latest Handoff Datetime" data for the two new rows and flip the
#    "Handoff Reason" from Ignored to Include

$wb = $excel.ActiveWorkbook

$newMdFile     = "f6f5a2a4-5f86-4251-8a44-9555b99eda99.md"
$newMdFile2    = "ffff856d4578-515f-42c7-8725-51aef913c338.md"
$configFile    = ".localization-config"
$readyStatus   = "Ready for handoff"
$notLocalized  = "Not to be localized"
$includeText   = "Include"
$epoch         = "0001-01-01 00:00:00"

$xlfZh   = "f6f5a2a4-5f86-4251-8a44-9555b99eda99.cab98be5212dacce1ff946b56f688b750c49bf20.zh-cn.xlf"
$xlfDe   = "f6f5a2a4-5f86-4251-8a44-9555b99eda99.cab98be5212dacce1ff946b56f688b750c49bf20.de-de.xlf"
$handoffZh = "2016-02-15 04:12:05"
$handoffDe = "2016-02-15 04:12:19"

$baseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/f6f5a2a4-5f86-4251-8a44-9555b99eda99"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# make room for the new file row; old row 3 (.localization-config) becomes row 4
$ws.Rows(3).Insert()

# clear every hyperlink on the sheet so we can rebuild them cleanly and in order
$ws.Range("A1").Hyperlinks.Delete()

# row 2 - renamed file, now ready for handoff
$ws.Range("B2").Value = $readyStatus
$ws.Range("C2").Value = $readyStatus

# row 3 - new file, also ready for handoff
$ws.Range("B3").Value = $readyStatus
$ws.Range("C3").Value = $readyStatus

# row 4 keeps its old values (.localization-config / Not to be localized),
# nothing else to change there

# rebuild hyperlinks in display order
$ws.Hyperlinks.Add($ws.Range("A2"), "$baseUrl/e2e/$newMdFile", "", "", $newMdFile)
$ws.Hyperlinks.Add($ws.Range("A3"), "$baseUrl/e2e/$newMdFile2", "", "", $newMdFile2)
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f9611c7f8fde1b0f51d8e99da3bc02453a1ad0c9/.localization-config", "", "", $configFile)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows(3).Insert()
$ws.Range("A1").Hyperlinks.Delete()

# row 2 - renamed file, handoff info populated, reason flips to Include
$ws.Range("B2").Value = $readyStatus
$ws.Range("D2").Value = $handoffZh
$ws.Range("G2").Value = $epoch
$ws.Range("H2").Value = $includeText

# row 3 - new file, same shape as row 2
$ws.Range("B3").Value = $readyStatus
$ws.Range("D3").Value = $handoffZh
$ws.Range("G3").Value = $epoch
$ws.Range("H3").Value = $includeText

# row 4 keeps its old values (.localization-config bookkeeping row)

$ws.Hyperlinks.Add($ws.Range("A2"), "$baseUrl/e2e/$newMdFile", "", "", $newMdFile)
$ws.Hyperlinks.Add($ws.Range("C2"), "$baseUrl/$xlfZh", "", "", $xlfZh)
$ws.Hyperlinks.Add($ws.Range("A3"), "$baseUrl/e2e/$newMdFile2", "", "", $newMdFile2)
$ws.Hyperlinks.Add($ws.Range("C3"), "$baseUrl/$xlfZh", "", "", $xlfZh)
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f9611c7f8fde1b0f51d8e99da3bc02453a1ad0c9/.localization-config", "", "", $configFile)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows(3).Insert()
$ws.Range("A1").Hyperlinks.Delete()

# row 2 - renamed file, handoff info populated, reason flips to Include
$ws.Range("B2").Value = $readyStatus
$ws.Range("D2").Value = $handoffDe
$ws.Range("G2").Value = $epoch
$ws.Range("H2").Value = $includeText

# row 3 - new file, same shape as row 2
$ws.Range("B3").Value = $readyStatus
$ws.Range("D3").Value = $handoffDe
$ws.Range("G3").Value = $epoch
$ws.Range("H3").Value = $includeText

# row 4 keeps its old values (.localization-config bookkeeping row)

$ws.Hyperlinks.Add($ws.Range("A2"), "$baseUrl/e2e/$newMdFile", "", "", $newMdFile)
$ws.Hyperlinks.Add($ws.Range("C2"), "$baseUrl/$xlfDe", "", "", $xlfDe)
$ws.Hyperlinks.Add($ws.Range("A3"), "$baseUrl/e2e/$newMdFile2", "", "", $newMdFile2)
$ws.Hyperlinks.Add($ws.Range("C3"), "$baseUrl/$xlfDe", "", "", $xlfDe)
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f9611c7f8fde1b0f51d8e99da3bc02453a1ad0c9/.localization-config", "", "", $configFile)

Write-Host "Report regenerated for handoff"
